# Revert "ograniczenie zasiegu smyrania": restore two list entries that had
# been removed, mark them done with a bold "[ZROBIONE]" prefix, and add two
# new "do ustalenia" items.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper RGB ints, computed from theme "Accent 6" (70AD47) tinted the
# same way Excel's built-in "40% - Accent6" swatch / font tint does ---
$fillColor = 11854021   # 0xC5E0B4 - Accent6, Lighter 60%
$fontColor = 3506772    # 0x548235 - Accent6, Darker 25%

# ----------------------------------------------------------------------
# E2: "fizyka, ..." gains a bold "[ZROBIONE]" prefix + green done styling
# ----------------------------------------------------------------------
$prefix = "[ZROBIONE]"
$fizykaText = "fizyka, kiedy działa na gracza nie przesuwa kamery. Śledzenie gracza trzeba uzależnić od jego pozycji, a nie od eventu poruszania (chyba)"
$e2text = $prefix + $fizykaText
$ws.Range("E2").Value = $e2text
$ws.Range("E2").Interior.Color = $fillColor
$ws.Range("E2").Font.Color = $fontColor
$e2bold = $ws.Range("E2").Characters(1, $prefix.Length)
$e2bold.Font.Bold = $true
$e2bold.Font.Color = $fontColor
$e2rest = $ws.Range("E2").Characters($prefix.Length + 1, $fizykaText.Length)
$e2rest.Font.Bold = $false
$e2rest.Font.Color = $fontColor

# ----------------------------------------------------------------------
# E3: "sterowanie, ..." gains a bold "[ZROBIONE]" prefix + green done styling
# ----------------------------------------------------------------------
$sterowanieText = "sterowanie, hover nad planszą zaznacza wszystkie pola nad jakimi jest myszka"
$e3text = $prefix + $sterowanieText
$ws.Range("E3").Value = $e3text
$ws.Range("E3").Interior.Color = $fillColor
$ws.Range("E3").Font.Color = $fontColor
$e3bold = $ws.Range("E3").Characters(1, $prefix.Length)
$e3bold.Font.Bold = $true
$e3bold.Font.Color = $fontColor
$e3rest = $ws.Range("E3").Characters($prefix.Length + 1, $sterowanieText.Length)
$e3rest.Font.Bold = $false
$e3rest.Font.Color = $fontColor

# ----------------------------------------------------------------------
# New "do ustalenia" items in column A (same look as A2)
# ----------------------------------------------------------------------
$ws.Range("A3").Value = "ekrany ładowania"
$ws.Range("A3").Style = "20% - akcent 3"

$ws.Range("A4").Value = "generowanie świata"
$ws.Range("A4").Style = "20% - akcent 3"

# ----------------------------------------------------------------------
# Row 2 grows taller to fit the now-longer wrapped text in E2
# ----------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 86.4

# ----------------------------------------------------------------------
# Selection moves to J4
# ----------------------------------------------------------------------
$ws.Range("J4").Select()
